# Generate Report for handoff
# The handoff of 44f9dfa4-... failed to transform; a new source file
# f939a530-8a03-4901-891d-bcc658750a13.md takes its place and the
# per-language status rows are reset to reflect the failed handoff
# (no handoff file/date yet, and the dependency is now Ignored).

$wb = $excel.ActiveWorkbook

$oldFile = "44f9dfa4-9b64-4e46-ae02-f2a609207392.md"
$newFile = "f939a530-8a03-4901-891d-bcc658750a13.md"
$oldStatus = "Ready for handoff"
$newStatus = "Handoff transform failed"
$resetDate = "0001-01-01 00:00:00"
$newReason = "Ignored"

function Update-FileNameHyperlink($ws, $cellAddr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $cellAddr) {
            $hl.TextToDisplay = $newFile
        }
    }
}

function Remove-HyperlinkAt($ws, $cellAddr) {
    $toDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $cellAddr) {
            $toDelete += $hl
        }
    }
    foreach ($hl in $toDelete) {
        $hl.Delete()
    }
}

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFile
Update-FileNameHyperlink $wsOverview '$A$2'
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# ---- Sheet "zh-cn" ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newFile
Update-FileNameHyperlink $wsZh '$A$2'
$wsZh.Range("B2").Value = $newStatus
# Latest Handoff File cell + its hyperlink are cleared (handoff never completed)
$wsZh.Range("C2").Clear()
Remove-HyperlinkAt $wsZh '$C$2'
# Latest Handoff Datetime resets
$wsZh.Range("D2").Value = $resetDate
# Handoff Reason is now "Ignored"
$wsZh.Range("H2").Value = $newReason

# ---- Sheet "de-de" ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newFile
Update-FileNameHyperlink $wsDe '$A$2'
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("C2").Clear()
Remove-HyperlinkAt $wsDe '$C$2'
$wsDe.Range("D2").Value = $resetDate
$wsDe.Range("H2").Value = $newReason
